$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-write B3 as a single, uniformly-formatted string (collapses the old
# two-run rich text "zota galax " + "usado seminovo" into one plain run).
$ws.Range("B3").Value = "zota galax usado seminovo"

# Update the minimum price for the rtx 3060 search.
$ws.Range("C3").Value = 2300

# Leave the selection where the loop finished (next row down from the table).
$ws.Range("C4").Select()
